$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared strings must be interned in this order to match the target:
#   127.0.0.1, ProxyServer_1, 000105001
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "ProxyServer_1"
$ws.Range("C2").Value = "ProxyServer_1"
$ws.Range("B2").Value = "000105001"

# Make sure text-like cells keep/gain the text number format (style index 1)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"

$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 5001

$ws.Range("G4").Select()
